# "Elimina antiguos EC y agrega nuevos y modifica Antigua BD"
# The old "Estado de Cuenta" period 2508 is retired and replaced by the new
# period 2509 for every worker row that still references it, and the
# "Periodo Mora" column in the data table is reformatted to be centered.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# 1) Update the "Periodo Mora" value 2508 -> 2509 wherever it appears.
$found = $ws.Cells.Find("2508")
if ($found -ne $null) {
    $firstAddress = $found.Address()
    while ($true) {
        $found.Value = "2509"
        $found = $ws.Cells.FindNext($found)
        if ($found -eq $null -or $found.Address() -eq $firstAddress) { break }
    }
}

# 2) Center-align the "Periodo Mora" column values for the data rows (E16:E22).
$ws.Range("E16:E22").HorizontalAlignment = -4108  # xlCenter
